$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Year" column header
$ws.Range("I1").Value = "Year"

# Add Year values for the two data rows
$ws.Range("I2").Value = 2011
$ws.Range("I3").Value = 2011

# Give I2 an explicit "General" number format (applies a numFmtId=0 explicit xf)
$ws.Range("I2").NumberFormat = "General"

# Clear the previous cell selection/active-cell on the sheet view
$ws.Range("A1").Select()
